$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.426422666666666
$ws.Range("H2").Value = 25.279268
$ws.Range("I2").Value = 0.1151758588783328
$ws.Range("J2").Value = 0.1151758588783328
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.864463
$ws.Range("N2").Value = 17.593389
$ws.Range("O2").Value = 0.6069167733108516
$ws.Range("P2").Value = 0.6069167733108515
$ws.Range("Q2").Value = 49.41644395102799
$ws.Range("R2").Value = 444.7479955592519
$ws.Range("S2").Value = 0.06990216063374373
$ws.Range("T2").Value = 0.06990216063374372

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.426422666666666
$ws.Range("H3").Value = 25.279268
$ws.Range("I3").Value = 0.1151758588783328
$ws.Range("J3").Value = 0.1151758588783328
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.691504
$ws.Range("N3").Value = 8.074512
$ws.Range("O3").Value = 0.278545354115671
$ws.Range("P3").Value = 0.278545354115671
$ws.Range("Q3").Value = 22.679750313024
$ws.Range("R3").Value = 204.117752817216
$ws.Range("S3").Value = 0.03208170039684176
$ws.Range("T3").Value = 0.03208170039684176

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.426422666666666
$ws.Range("H4").Value = 25.279268
$ws.Range("I4").Value = 0.1151758588783328
$ws.Range("J4").Value = 0.1151758588783328
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.106746666666667
$ws.Range("N4").Value = 3.32024
$ws.Range("O4").Value = 0.1145378725734776
$ws.Range("P4").Value = 0.1145378725734776
$ws.Range("Q4").Value = 9.325915198257777
$ws.Range("R4").Value = 83.93323678432
$ws.Range("S4").Value = 0.01319199784774732
$ws.Range("T4").Value = 0.01319199784774731

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.78712033333333
$ws.Range("H5").Value = 35.361361
$ws.Range("I5").Value = 0.1611112760180311
$ws.Range("J5").Value = 0.1611112760180311
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.864463
$ws.Range("N5").Value = 17.593389
$ws.Range("O5").Value = 0.6069167733108516
$ws.Range("P5").Value = 0.6069167733108515
$ws.Range("Q5").Value = 69.12513107138101
$ws.Range("R5").Value = 622.126179642429
$ws.Range("S5").Value = 0.09778113578485743
$ws.Range("T5").Value = 0.09778113578485742

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.78712033333333
$ws.Range("H6").Value = 35.361361
$ws.Range("I6").Value = 0.1611112760180311
$ws.Range("J6").Value = 0.1611112760180311
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.691504
$ws.Range("N6").Value = 8.074512
$ws.Range("O6").Value = 0.278545354115671
$ws.Range("P6").Value = 0.278545354115671
$ws.Range("Q6").Value = 31.725081525648
$ws.Range("R6").Value = 285.525733730832
$ws.Range("S6").Value = 0.04487679743047009
$ws.Range("T6").Value = 0.04487679743047009

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.78712033333333
$ws.Range("H7").Value = 35.361361
$ws.Range("I7").Value = 0.1611112760180311
$ws.Range("J7").Value = 0.1611112760180311
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.106746666666667
$ws.Range("N7").Value = 3.32024
$ws.Range("O7").Value = 0.1145378725734776
$ws.Range("P7").Value = 0.1145378725734776
$ws.Range("Q7").Value = 13.04535613851556
$ws.Range("R7").Value = 117.40820524664
$ws.Range("S7").Value = 0.01845334280270362
$ws.Range("T7").Value = 0.01845334280270362

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 52.94781866666667
$ws.Range("H8").Value = 158.843456
$ws.Range("I8").Value = 0.7237128651036362
$ws.Range("J8").Value = 0.7237128651036362
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.864463
$ws.Range("N8").Value = 17.593389
$ws.Range("O8").Value = 0.6069167733108516
$ws.Range("P8").Value = 0.6069167733108515
$ws.Range("Q8").Value = 310.510523501376
$ws.Range("R8").Value = 2794.594711512384
$ws.Range("S8").Value = 0.4392334768922505
$ws.Range("T8").Value = 0.4392334768922504

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 52.94781866666667
$ws.Range("H9").Value = 158.843456
$ws.Range("I9").Value = 0.7237128651036362
$ws.Range("J9").Value = 0.7237128651036362
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.691504
$ws.Range("N9").Value = 8.074512
$ws.Range("O9").Value = 0.278545354115671
$ws.Range("P9").Value = 0.278545354115671
$ws.Range("Q9").Value = 142.509265732608
$ws.Range("R9").Value = 1282.583391593472
$ws.Range("S9").Value = 0.2015868562883592
$ws.Range("T9").Value = 0.2015868562883592

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 52.94781866666667
$ws.Range("H10").Value = 158.843456
$ws.Range("I10").Value = 0.7237128651036362
$ws.Range("J10").Value = 0.7237128651036362
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.106746666666667
$ws.Range("N10").Value = 3.32024
$ws.Range("O10").Value = 0.1145378725734776
$ws.Range("P10").Value = 0.1145378725734776
$ws.Range("Q10").Value = 58.59982181660445
$ws.Range("R10").Value = 527.3983963494401
$ws.Range("S10").Value = 0.08289253192302665
$ws.Range("T10").Value = 0.08289253192302663
